$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any existing content first to avoid leftovers from the old smaller range
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 6565.833333333333, 7112, 6140, 0.1588538487752279),
    @(1, 6427.966666666666, 6829, 5647, 0.1592517852783203),
    @(2, 6250.4,            6755, 5387, 0.1611189206441243),
    @(3, 6833.4,            7438, 6475, 0.1606688102086385),
    @(4, 6334.6,            6855, 5532, 0.1657414436340332),
    @(5, 6318.966666666666, 6790, 5859, 0.1602704763412476),
    @(6, 6678.166666666667, 7159, 5812, 0.1639129161834717),
    @(7, 6623.933333333333, 7157, 5969, 0.1604065736134847),
    @(8, 6418.9,            6988, 5754, 0.1607615391413371),
    @(9, 6653.5,            7311, 5966, 0.160109003384908)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
